# Generate Report for Archive
#
# The localization job for "b5b58f96-a234-43cf-ba44-6d640241a84c.md" moved
# out of the handoff stage and back into translation, so its Status
# changes from "Ready for handoff" to "In Translation" everywhere it is
# reported: the per-language tables (zh-cn, de-de) and the roll-up
# Overview table.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 4 is the b5b58f96-... file, column C is "Status"
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C4").Value = "In Translation"

# de-de sheet: row 4 is the b5b58f96-... file, column C is "Status"
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C4").Value = "In Translation"

# Overview sheet: row 4 is the b5b58f96-... file, columns E (zh-cn) and
# F (de-de) mirror the per-language status
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
